$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.582.28'
$ws.Range('E2').Value = '  -2.05%  '

# Row 3
$ws.Range('D3').Value = '2.997.86'
$ws.Range('E3').Value = '  -1.07%  '

# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').Value = '''595.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.49%  '

# Row 6
$ws.Range('D6').Value = '''144.25'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.09%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('D8').Value = '''0.522'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.83%  '

# Row 9
$ws.Range('D9').Value = '2.997.36'
$ws.Range('E9').Value = '  -1.03%  '

# Row 10
$ws.Range('E10').Value = '  -2.12%  '

# Row 11
$ws.Range('D11').Value = '''5.91'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.63%  '

# Row 12
$ws.Range('D12').Value = '''0.461'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.85%  '

# Row 13
$ws.Range('D13').Value = '''0.0000229'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.27%  '

# Row 14
$ws.Range('D14').Value = '''34.31'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.86%  '

# Row 15
$ws.Range('E15').Value = '  +2.09%  '

# Row 16
$ws.Range('D16').Value = '3.494.78'
$ws.Range('E16').Value = '  -1.15%  '

# Row 17
$ws.Range('D17').Value = '''7.05'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.34%  '

# Row 18
$ws.Range('D18').Value = '61.600.72'
$ws.Range('E18').Value = '  -1.98%  '

# Row 19
$ws.Range('D19').Value = '3.002.21'
$ws.Range('E19').Value = '  -1.04%  '

# Row 20
$ws.Range('D20').Value = '''453.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.37%  '

# Row 21
$ws.Range('D21').Value = '''14.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.19%  '

# Row 22
$ws.Range('D22').Value = '''0.687'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.74%  '

# Row 23
$ws.Range('D23').Value = '''7.36'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.66%  '

# Row 24
$ws.Range('D24').Value = '''82.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.43%  '

# Row 25
$ws.Range('D25').Value = '''2.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.28%  '

# Row 26
$ws.Range('D26').Value = '''10.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.03%  '

# Row 27
$ws.Range('D27').Value = '''12.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.18%  '

# Row 28
$ws.Range('E28').Value = '  +0.03%  '

# Row 29
$ws.Range('D29').Value = '''2.67'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.32%  '

# Row 30
$ws.Range('D30').Value = '''1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.04%  '

# Row 31
$ws.Range('D31').Value = '''7.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.80%  '

# Row 32
$ws.Range('D32').Value = '''2.08'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.33%  '

# Row 33
$ws.Range('D33').Value = '''27.53'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.04%  '

# Row 35
$ws.Range('D35').Value = '0.0₃0835'
$ws.Range('E35').Value = '  +3.63%  '

# Row 36
$ws.Range('E36').Value = '  -2.05%  '

# Row 37
$ws.Range('D37').Value = '''5.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.24%  '

# Row 38
$ws.Range('D38').Value = '''9.26'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.64%  '

# Row 39
$ws.Range('D39').Value = '''2.07'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.80%  '

# Row 40
$ws.Range('D40').Value = '''50.27'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.18%  '

# Row 41
$ws.Range('D41').Value = '''2.90'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.17%  '

# Row 42
$ws.Range('E42').Value = '  +7.60%  '

# Row 43
$ws.Range('D43').Value = '''395.20'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.46%  '

# Row 44
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').Value = '''39.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.56%  '

# Row 45
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0354'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.61%  '

# Row 46
$ws.Range('E46').Value = '  -3.91%  '

# Row 47
$ws.Range('D47').Value = '2.722.01'
$ws.Range('E47').Value = '  -2.96%  '

# Row 48
$ws.Range('D48').Value = '''133.27'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.41%  '

# Row 49
$ws.Range('E49').Value = '  +0.13%  '

# Row 50
$ws.Range('D50').Value = '''0.107'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.81%  '

# Row 51
$ws.Range('D51').Value = '''2.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.24%  '
